# "added Inori and Marco"
#
# On the "Characters" sheet, fill in the stat columns (Health, MP, Attack,
# Defense, Resistance, Skill, Speed -> columns G:M) for two characters that
# previously had no stats entered:
#   - Inori (row 31)
#   - Marco (row 35)
# The Total column (N) already holds a shared SUM(Gn:Mn) formula for every
# row, so it recalculates on its own once G:M are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Characters")

# Inori - row 31: Health, MP, Attack, Defense, Resistance, Skill, Speed
$ws.Range("G31").Value = 60
$ws.Range("H31").Value = 80
$ws.Range("I31").Value = 55
$ws.Range("J31").Value = 45
$ws.Range("K31").Value = 45
$ws.Range("L31").Value = 40
$ws.Range("M31").Value = 65

# Marco - row 35: Health, MP, Attack, Defense, Resistance, Skill, Speed
$ws.Range("G35").Value = 70
$ws.Range("H35").Value = 50
$ws.Range("I35").Value = 65
$ws.Range("J35").Value = 50
$ws.Range("K35").Value = 70
$ws.Range("L35").Value = 40
$ws.Range("M35").Value = 40

# Both characters now have art, so flip their "Has Art" form-control checkbox.
$ws.Shapes.Item("Check Box 67").ControlFormat.Value = 1
$ws.Shapes.Item("Check Box 71").ControlFormat.Value = 1

# Leave the view where editing ended up: scrolled down with G36 selected.
$null = $ws.Range("G36").Select()
$excel.ActiveWindow.ScrollRow = 22
